# SIS-196 — adiciona colunas do Comitê ao dashboard
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (L:M) right after the existing "Boletos Pagos" total
# column (K). Excel shifts everything from L onward two columns to the right
# and carries the style (s="3", bold/centered header) of the pushed-aside
# cells onto the freshly inserted ones, matching the original author's edit.
$ws.Columns("L:M").Insert()

# Fix the typo in the (now) K1 header and populate the two new header cells.
$ws.Range("K1").Value = "Valor Total de Contratos com Boeltos Pagos"
$ws.Range("L1").Value = "Aprovadas Comite"
$ws.Range("M1").Value = "Valor Total Comite"

# Restore the view state the author left the sheet in.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.Zoom = 85
$ws.Range("N2").Select()
